# Updates the crypto price/volume snapshot on Sheet1 (cols D "Price" and
# E "Volume(1h)") to the refreshed values from the latest data pull.
#
# The source cells are plain text (inline strings) even though most of
# them look numeric (e.g. "49.70", "0.01080") -- the trailing zeros and
# exact decimal representation matter, so a bare `.Value = "..."` assignment
# would let Excel auto-coerce them into real floating point numbers and
# lose precision/formatting. To keep them as text without altering the
# cell's style, we temporarily force a text number format, assign the
# value, then restore the style to "Normal" (which matches the workbook's
# original default style for these cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "319.87"
Set-TextValue "D3" "49.08"
Set-TextValue "E3" "10.75%"
Set-TextValue "D4" "5.320"
Set-TextValue "E4" "3.87%"
Set-TextValue "D5" "0.07987"
Set-TextValue "E5" "0.77%"
Set-TextValue "D6" "4.601"
Set-TextValue "E6" "3.62%"
Set-TextValue "D7" "1.338"
Set-TextValue "E7" "26.60%"
Set-TextValue "E8" "3.35%"
Set-TextValue "D9" "0.1248"
Set-TextValue "E9" "-1.58%"
Set-TextValue "D10" "0.1974"
Set-TextValue "E10" "5.31%"
Set-TextValue "D11" "0.09647"
Set-TextValue "E11" "5.39%"
Set-TextValue "D12" "0.04532"
Set-TextValue "E12" "8.99%"
Set-TextValue "E13" "0.19%"
Set-TextValue "D14" "0.001307"
Set-TextValue "E14" "0.58%"
Set-TextValue "D15" "0.04219"
Set-TextValue "E15" "1.06%"
Set-TextValue "D16" "0.005831"
Set-TextValue "E16" "1.00%"
Set-TextValue "D17" "3.348"
Set-TextValue "E17" "-1.06%"
Set-TextValue "E18" "5.66%"
Set-TextValue "D19" "0.3471"
Set-TextValue "E19" "1.55%"
Set-TextValue "D20" "8.089"
Set-TextValue "E20" "0.88%"
Set-TextValue "D21" "0.1402"
Set-TextValue "E21" "2.00%"
Set-TextValue "D22" "0.3006"
Set-TextValue "E22" "7.60%"
Set-TextValue "D23" "0.001296"
Set-TextValue "E23" "2.13%"
Set-TextValue "D24" "0.004232"
Set-TextValue "E24" "-6.17%"
Set-TextValue "D25" "0.0001354"
Set-TextValue "E25" "1.22%"
Set-TextValue "D26" "0.0003548"
Set-TextValue "E26" "-95.21%"
Set-TextValue "D38" "0.02656"
Set-TextValue "E38" "0.33%"
Set-TextValue "D39" "0.05946"
Set-TextValue "E39" "10.94%"
Set-TextValue "D40" "0.01079"
Set-TextValue "E40" "94.56%"
Set-TextValue "D41" "0.008043"
Set-TextValue "E41" "4.25%"
Set-TextValue "D42" "0.1461"
Set-TextValue "E42" "5.71%"
Set-TextValue "D43" "0.007536"
Set-TextValue "E43" "3.56%"
Set-TextValue "D44" "0.007959"
Set-TextValue "E44" "-3.86%"
Set-TextValue "D45" "0.3209"
Set-TextValue "E45" "6.04%"
Set-TextValue "D46" "0.00007026"
Set-TextValue "E46" "5.45%"
Set-TextValue "D47" "0.00000000752"
Set-TextValue "E47" "1.26%"
Set-TextValue "D48" "0.05592"
Set-TextValue "E48" "-9.44%"
Set-TextValue "D49" "0.004010"
Set-TextValue "E49" "1.22%"
Set-TextValue "D50" "0.00002106"
Set-TextValue "E50" "1.26%"
Set-TextValue "D51" "0.0002006"
Set-TextValue "E51" "1.26%"
